# "add user formular dialog"
# Adds two new toolbar-button rows ("Track Bolk Update" / "Curve Bolk Update")
# right after the existing "Auto Size Track" row, and renames the
# "AutoSizeTrackCheckbox" code to "AutoSizeTrackButton".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LogPlotTab")

# Insert two fresh rows directly below row 46 (the "Auto Size Track" button
# row). This pushes the old rows 47-50 down to 49-52 and lets the engine
# rewrite all the dependent formula references automatically.
$ws.Rows.Item(47).Insert()
$ws.Rows.Item(47).Insert()

# Rename the existing Auto Size Track row's code from a Checkbox-style name
# to a Button-style name (it has always been a wiButton).
$ws.Cells.Item(46, 2).Value = "AutoSizeTrackButton"

# New codes/labels, in the same order the author's form must have written
# them (code column first, then image names / tooltip labels).
$ws.Cells.Item(47, 2).Value = "TrackBolkUpdateButton"
$ws.Cells.Item(48, 2).Value = "CurveBolkUpdateButton"
$ws.Cells.Item(47, 5).Value = "track_bolk_16x16"
$ws.Cells.Item(47, 11).Value = "Track Bolk Update"
$ws.Cells.Item(48, 11).Value = "Curve Bolk Update"
$ws.Cells.Item(48, 5).Value = "curve_bolk_16x16"

# Remaining cells for the two new rows.
$ws.Cells.Item(47, 1).Value = 7.2
$ws.Cells.Item(47, 3).Value = "wiButton"
$ws.Cells.Item(47, 4).Formula = '=REPLACE(C47, 1, 2, "")'
$ws.Cells.Item(47, 6).Formula = '=SUBSTITUTE(E47,"_","-")'
$ws.Cells.Item(47, 10).Value = "small"

$ws.Cells.Item(48, 1).Value = 7.3
$ws.Cells.Item(48, 3).Value = "wiButton"
$ws.Cells.Item(48, 4).Formula = '=REPLACE(C48, 1, 2, "")'
$ws.Cells.Item(48, 6).Formula = '=SUBSTITUTE(E48,"_","-")'
$ws.Cells.Item(48, 10).Value = "small"

# Match the author's final selection.
$ws.Range("E48").Select()
